# Auto-generated edit script applying the Adamantoise_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for specific rows
# across the ALC, ARM, BSM, CRP, CUL, GSM, WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 545.4
$ws.Range("I58").Value = 545.4
$ws.Range("K58").Value = 1636.2
$ws.Range("M58").Value = -1486.2
$ws.Range("H61").Value = 130.3
$ws.Range("I61").Value = 138.75
$ws.Range("J61").Value = 96.5
$ws.Range("K61").Value = 416.25
$ws.Range("L61").Value = 289.5
$ws.Range("M61").Value = -244.25
$ws.Range("N61").Value = -633.5
$ws.Range("H99").Value = 1080.6666
$ws.Range("I99").Value = 1108.6666
$ws.Range("J99").Value = 1066.6666
$ws.Range("K99").Value = 3325.9998
$ws.Range("L99").Value = 3199.9998
$ws.Range("M99").Value = -1827.9998
$ws.Range("N99").Value = -6195.9998
$ws.Range("H100").Value = 2717.889
$ws.Range("I100").Value = 1840.1538
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 1840.1538
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -1299.1538
$ws.Range("N100").Value = -6082
$ws.Range("H107").Value = 3207.6365
$ws.Range("I107").Value = 3207.6365
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3207.6365
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1287.6365
$ws.Range("H115").Value = 924.8570999999999
$ws.Range("I115").Value = 743.75
$ws.Range("J115").Value = 1166.3334
$ws.Range("K115").Value = 2231.25
$ws.Range("L115").Value = 3499.0002
$ws.Range("M115").Value = -664.25
$ws.Range("N115").Value = -6633.0002
$ws.Range("H116").Value = 13179.5625
$ws.Range("I116").Value = 16198.25
$ws.Range("J116").Value = 4123.5
$ws.Range("K116").Value = 16198.25
$ws.Range("L116").Value = 4123.5
$ws.Range("M116").Value = -12756.25
$ws.Range("N116").Value = -11007.5
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("H132").Value = 1579.9166
$ws.Range("I132").Value = 1431.4348
$ws.Range("K132").Value = 4294.3044
$ws.Range("M132").Value = -1764.3044
$ws.Range("H138").Value = 3087.12
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 3087.12
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 9261.360000000001
$ws.Range("N138").Value = -19541.36
$ws.Range("N107").ClearContents()
$ws.Range("N129").ClearContents()
$ws.Range("M138").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 69170.5
$ws.Range("I31").Value = 39997.5
$ws.Range("J31").Value = 83757
$ws.Range("K31").Value = 39997.5
$ws.Range("L31").Value = 83757
$ws.Range("M31").Value = -39703.5
$ws.Range("N31").Value = -84345
$ws.Range("H32").Value = 18519556
$ws.Range("I32").Value = 20409010
$ws.Range("K32").Value = 20409010
$ws.Range("M32").Value = -20408723
$ws.Range("H102").Value = 2719.6
$ws.Range("I102").Value = 2032.6666
$ws.Range("J102").Value = 3750
$ws.Range("K102").Value = 2032.6666
$ws.Range("L102").Value = 3750
$ws.Range("M102").Value = -410.6666
$ws.Range("N102").Value = -6994
$ws.Range("H110").Value = 2332.2173
$ws.Range("I110").Value = 1620.4
$ws.Range("K110").Value = 1620.4
$ws.Range("M110").Value = 424.5999999999999
$ws.Range("H126").Value = 9999
$ws.Range("I126").Value = 9999
$ws.Range("K126").Value = 29997
$ws.Range("M126").Value = -27527
$ws.Range("H130").Value = 49403
$ws.Range("J130").Value = 49403
$ws.Range("L130").Value = 49403
$ws.Range("N130").Value = -59443

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2108.6667
$ws.Range("I105").Value = 1821.3334
$ws.Range("J105").Value = 3832.6667
$ws.Range("K105").Value = 1821.3334
$ws.Range("L105").Value = 3832.6667
$ws.Range("M105").Value = -74.33339999999998
$ws.Range("N105").Value = -7326.6667
$ws.Range("H107").Value = 3587.85
$ws.Range("I107").Value = 3575.8235
$ws.Range("J107").Value = 3656
$ws.Range("K107").Value = 3575.8235
$ws.Range("L107").Value = 3656
$ws.Range("M107").Value = -1655.8235
$ws.Range("N107").Value = -7496
$ws.Range("H134").Value = 3178.9524
$ws.Range("I134").Value = 3171.4736
$ws.Range("K134").Value = 9514.4208
$ws.Range("M134").Value = -6979.4208

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1914
$ws.Range("I19").Value = 936.6667
$ws.Range("J19").Value = 4008.2856
$ws.Range("K19").Value = 936.6667
$ws.Range("L19").Value = 4008.2856
$ws.Range("M19").Value = -766.6667
$ws.Range("N19").Value = -4348.2856
$ws.Range("H24").Value = 1914
$ws.Range("I24").Value = 936.6667
$ws.Range("J24").Value = 4008.2856
$ws.Range("K24").Value = 936.6667
$ws.Range("L24").Value = 4008.2856
$ws.Range("M24").Value = -766.6667
$ws.Range("N24").Value = -4348.2856
$ws.Range("H31").Value = 4605
$ws.Range("I31").Value = 2251.5
$ws.Range("J31").Value = 5487.5625
$ws.Range("K31").Value = 2251.5
$ws.Range("L31").Value = 5487.5625
$ws.Range("M31").Value = -1956.5
$ws.Range("N31").Value = -6077.5625
$ws.Range("H34").Value = 4605
$ws.Range("I34").Value = 2251.5
$ws.Range("J34").Value = 5487.5625
$ws.Range("K34").Value = 2251.5
$ws.Range("L34").Value = 5487.5625
$ws.Range("M34").Value = -2049.5
$ws.Range("N34").Value = -5891.5625
$ws.Range("H68").Value = 69955.60000000001
$ws.Range("J68").Value = 69955.60000000001
$ws.Range("L68").Value = 69955.60000000001
$ws.Range("N68").Value = -71453.60000000001
$ws.Range("H71").Value = 69955.60000000001
$ws.Range("J71").Value = 69955.60000000001
$ws.Range("L71").Value = 209866.8
$ws.Range("N71").Value = -217354.8
$ws.Range("H94").Value = 1585
$ws.Range("I94").Value = 816.6667
$ws.Range("K94").Value = 816.6667
$ws.Range("M94").Value = -365.6667
$ws.Range("H122").Value = 6759.5454
$ws.Range("I122").Value = 6995.5
$ws.Range("K122").Value = 20986.5
$ws.Range("M122").Value = -18536.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 70.60869599999999
$ws.Range("I38").Value = 65.1579
$ws.Range("J38").Value = 96.5
$ws.Range("K38").Value = 195.4737
$ws.Range("L38").Value = 289.5
$ws.Range("M38").Value = 151.5263
$ws.Range("N38").Value = -983.5
$ws.Range("H47").Value = 91000000
$ws.Range("I47").Value = 91000000
$ws.Range("K47").Value = 273000000
$ws.Range("M47").Value = -272999569
$ws.Range("H107").Value = 680.0526
$ws.Range("J107").Value = 699.25
$ws.Range("L107").Value = 2097.75
$ws.Range("N107").Value = -5937.75
$ws.Range("H131").Value = 1578.4314
$ws.Range("J131").Value = 1680.1428
$ws.Range("L131").Value = 5040.428400000001
$ws.Range("N131").Value = -15120.4284
$ws.Range("H133").Value = 4406.6
$ws.Range("I133").Value = 4190.8184
$ws.Range("K133").Value = 12572.4552
$ws.Range("M133").Value = -7512.4552
$ws.Range("H137").Value = 11551.7
$ws.Range("J137").Value = 100000
$ws.Range("L137").Value = 300000
$ws.Range("N137").Value = -310200
$ws.Range("H139").Value = 3679.3635
$ws.Range("I139").Value = 3347.3
$ws.Range("K139").Value = 10041.9
$ws.Range("M139").Value = -4901.900000000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1287.1111
$ws.Range("I97").Value = 512.1429000000001
$ws.Range("K97").Value = 512.1429000000001
$ws.Range("M97").Value = -16.14290000000005
$ws.Range("H102").Value = 1474.2
$ws.Range("I102").Value = 1364.56
$ws.Range("K102").Value = 1364.56
$ws.Range("M102").Value = 257.4400000000001
$ws.Range("H134").Value = 94108.336
$ws.Range("J134").Value = 94108.336
$ws.Range("L134").Value = 282325.008
$ws.Range("N134").Value = -287395.008

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 877.1579
$ws.Range("I100").Value = 553.2727
$ws.Range("J100").Value = 1322.5
$ws.Range("K100").Value = 1106.5454
$ws.Range("L100").Value = 2645
$ws.Range("M100").Value = -565.5454
$ws.Range("N100").Value = -3727
$ws.Range("H136").Value = 2628.4644
$ws.Range("I136").Value = 1879.4
$ws.Range("J136").Value = 3044.611
$ws.Range("K136").Value = 5638.200000000001
$ws.Range("L136").Value = 9133.832999999999
$ws.Range("M136").Value = -3088.200000000001
$ws.Range("N136").Value = -14233.833

Write-Host "Applied Adamantoise_Profits updates across ALC, ARM, BSM, CRP, CUL, GSM, WVR sheets"